# Cofina contrat_pep_individual_business.docx - V1 cofina contrat 06 dec 2024
#
# This script rewrites the "est representee par ..." paragraph so that the
# placeholders switch from the `individual_business.*` family to the new
# `verbal_trial.*` / `representative_*` family, and adds a number of new
# fields (birth date/place, identity document details, office delivery,
# home address, phone number) in the representative's legal description.

$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Locate the block to rewrite: from "${individual_business.civility}"
#    through to the end of the paragraph ("..., gerant ayant pleins
#    pouvoirs a l'effet des presentes,"). We keep "est representee par"
#    (and the bookmark around the final "r") untouched, then replace the
#    remainder wholesale with plain text, then re-apply bold formatting
#    to the placeholder spans.
# ----------------------------------------------------------------------

$startRange = $d.Content
$null = $startRange.Find.Execute("`${individual_business.civility`}")
$blockStart = $startRange.Start

$endRange = $d.Content
$null = $endRange.Find.Execute("gérant ayant pleins pouvoirs à l'effet des présentes")
$blockEnd = $endRange.End

$block = $d.Range($blockStart, $blockEnd)

$newText = "`${verbal_trial.civility} `${verbal_trial.applicant_last_name} `${verbal_trial.applicant_first_name} né le `${representative_birth_date}  à  `${representative_birth_place}, titulaire de `${representative_type_of_identity_document}  N°`${representative_number_of_identity_document} délivré le `${representative_date_of_issue_of_identity_document} par `${representative_office_delivery}, domicilié `${representative_home_address}, et répondant au `${representative_phone_number}, gérant ayant pleins pouvoirs à l'effet des présentes"

$block.Text = $newText

Write-Host "Replaced block text. Trailing char now: [$($d.Range($block.End, $block.End + 1).Text)]"

# ----------------------------------------------------------------------
# 2) Re-apply bold formatting to the placeholder spans, in document
#    order, by searching forward from the previous match each time.
# ----------------------------------------------------------------------

$boldSpans = @(
    " `${verbal_trial.applicant_last_name} `${verbal_trial.applicant_first_name} ",
    "`${representative_birth_date}  ",
    "`${representative_birth_place}, ",
    "`${representative_type_of_identity_document}  N°`${representative_number_of_identity_document} ",
    "`${representative_date_of_issue_of_identity_document} ",
    " `${representative_office_delivery}",
    "`${representative_home_address}, ",
    "`${representative_phone_number}, "
)

$searchFrom = $blockStart
foreach ($span in $boldSpans) {
    $searchRange = $d.Range($searchFrom, $d.Content.End)
    $found = $searchRange.Find.Execute($span)
    if (-not $found) {
        Write-Host "NOT FOUND: $span"
        continue
    }
    $searchRange.Bold = 1
    $searchFrom = $searchRange.End
}

Write-Host "Bold spans applied."
